$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
  @("ECs","Gnai2","P2ry12","ECs",3,1,151.7260716666667,455.178215,0.2700739458961593,0.2783366498663096,2,0.6666666666666666,0.1444173333333333,0.433252,0.003554110791691976,0.003583730627916089,21.91187466724222,197.20687200518,0.0009598727256643746,0.0009974835769974502),
  @("ECs","Gnai2","P2ry12","M1",3,1,151.7260716666667,455.178215,0.2700739458961593,0.2783366498663096,3,1,17.275143,51.825429,0.4251412953499726,0.4286844081785905,2621.089584869915,23589.80626382924,0.1148195871985716,0.1193185820223505),
  @("ECs","Gnai2","P2ry12","M2",3,1,151.7260716666667,455.178215,0.2700739458961593,0.2783366498663096,3,1,22.20679833333334,66.620395,0.5465093405599564,0.551063930473184,3369.350275410548,30324.15247869493,0.1475979340741354,0.153381288270067),
  @("ECs","Gnai2","P2ry12","sCs",3,1,151.7260716666667,455.178215,0.2700739458961593,0.2783366498663096,2,1,1.0075275,2.015055,0.02479525329837903,0.01666793072030932,152.8681896711375,917.2091380268249,0.006696551897787984,0.004639295996894641),
  @("FAPs","Gnai2","P2ry12","ECs",3,1,82.248871,246.746613,0.146403824289839,0.150882936320401,2,0.6666666666666666,0.1444173333333333,0.433252,0.003554110791691976,0.003583730627916089,11.87816261949733,106.903463575476,0.0005203354118534925,0.0005407238001213339),
  @("FAPs","Gnai2","P2ry12","M1",3,1,82.248871,246.746613,0.146403824289839,0.150882936320401,3,1,17.275143,51.825429,0.4251412953499726,0.4286844081785905,1420.861008113553,12787.74907302198,0.06224231150277192,0.06468116226075905),
  @("FAPs","Gnai2","P2ry12","M2",3,1,82.248871,246.746613,0.146403824289839,0.150882936320401,3,1,22.20679833333334,66.620395,0.5465093405599564,0.551063930473184,1826.484091441348,16438.35682297213,0.08001105746809561,0.0831461439300553),
  @("FAPs","Gnai2","P2ry12","sCs",3,1,82.248871,246.746613,0.146403824289839,0.150882936320401,2,1,1.0075275,2.015055,0.02479525329837903,0.01666793072030932,82.86799937645249,497.207996258715,0.003630119907117934,0.002514906329465287),
  @("M1","Gnai2","P2ry12","ECs",3,1,123.444321,370.332963,0.2197321429647646,0.2264546783208506,2,0.6666666666666666,0.1444173333333333,0.433252,0.003554110791691976,0.003583730627916089,17.827499653964,160.447496885676,0.0007809523805926739,0.0008115525665333179),
  @("M1","Gnai2","P2ry12","M1",3,1,123.444321,370.332963,0.2197321429647646,0.2264546783208506,3,1,17.275143,51.825429,0.4251412953499726,0.4286844081785905,2132.518297812903,19192.66468031613,0.09341720789006538,0.09707758975524694),
  @("M1","Gnai2","P2ry12","M2",3,1,123.444321,370.332963,0.2197321429647646,0.2264546783208506,3,1,22.20679833333334,66.620395,0.5465093405599564,0.551063930473184,2741.303141842265,24671.72827658039,0.1200856685514995,0.1247910051095285),
  @("M1","Gnai2","P2ry12","sCs",3,1,123.444321,370.332963,0.2197321429647646,0.2264546783208506,2,1,1.0075275,2.015055,0.02479525329837903,0.01666793072030932,124.3735481263275,746.241288757965,0.005448314142606972,0.003774530889541871),
  @("M2","Gnai2","P2ry12","ECs",3,1,154.3429766666667,463.02893,0.2747320633285943,0.2831372788071194,2,0.6666666666666666,0.1444173333333333,0.433252,0.003554110791691976,0.003583730627916089,22.28980110892888,200.60820998036,0.0009764281910999603,0.001014687737965891),
  @("M2","Gnai2","P2ry12","M1",3,1,154.3429766666667,463.02893,0.2747320633285943,0.2831372788071194,3,1,17.275143,51.825429,0.4251412953499726,0.4286844081785905,2666.29699296233,23996.67293666097,0.1167999452776893,0.1213765367987266),
  @("M2","Gnai2","P2ry12","M2",3,1,154.3429766666667,463.02893,0.2747320633285943,0.2831372788071194,3,1,22.20679833333334,66.620395,0.5465093405599564,0.551063930473184,3427.463357003039,30847.17021302735,0.1501436387603862,0.156026741722933),
  @("M2","Gnai2","P2ry12","sCs",3,1,154.3429766666667,463.02893,0.2747320633285943,0.2831372788071194,2,1,1.0075275,2.015055,0.02479525329837903,0.01666793072030932,155.504793423525,933.02876054115,0.006812051099418805,0.004719312547493972),
  @("sCs","Gnai2","P2ry12","ECs",2,1,50.0323125,100.064625,0.08905802352064279,0.06118845668531954,2,0.6666666666666666,0.1444173333333333,0.433252,0.003554110791691976,0.003583730627916089,7.22553315175,43.3531989105,0.0003165220824814744,0.0002192829462980966),
  @("sCs","Gnai2","P2ry12","M1",2,1,50.0323125,100.064625,0.08905802352064279,0.06118845668531954,3,1,17.275143,51.825429,0.4251412953499726,0.4286844081785905,864.3153530581875,5185.892118349126,0.0378622434808744,0.02623053734150753),
  @("sCs","Gnai2","P2ry12","M2",2,1,50.0323125,100.064625,0.08905802352064279,0.06118845668531954,3,1,22.20679833333334,66.620395,0.5465093405599564,0.551063930473184,1111.057473837813,6666.344843026875,0.04867104170583958,0.03371875144060036),
  @("sCs","Gnai2","P2ry12","sCs",2,1,50.0323125,100.064625,0.08905802352064279,0.06118845668531954,2,1,1.0075275,2.015055,0.02479525329837903,0.01666793072030932,50.40893073234375,201.635722929375,0.002208216251447336,0.001019884956913554)
)


$r = 2
foreach ($row in $rows) {
    $c = 1
    foreach ($val in $row) {
        $ws.Cells.Item($r, $c).Value = $val
        $c++
    }
    $r++
}
